$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header in G1 from "mida" to "tamaño"
$ws.Range("G1").Value = "tamaño"

# Update the active selection to G3 (as reflected in the saved view state)
$ws.Range("G3").Select()
